$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 613, shifting all subsequent rows (613-666) down by one.
$ws.Rows.Item(613).EntireRow.Insert()

# Populate the newly inserted row 613 with the new record's data.
$ws.Range("A613").Value2 = 7
$ws.Range("B613").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C613").Value2 = "Ñuble"
$ws.Range("D613").Value2 = 45013
$ws.Range("E613").Value2 = 16
$ws.Range("F613").Value2 = 100112020
$ws.Range("G613").Value2 = "Tomate"
$ws.Range("H613").Value2 = "Larga vida"
$ws.Range("I613").Value2 = "Primera"
$ws.Range("J613").Value2 = 50
$ws.Range("K613").Value2 = 7000
$ws.Range("L613").Value2 = 7000
$ws.Range("M613").Value2 = 7000
$ws.Range("N613").Value2 = "$/bandeja 18 kilos"
$ws.Range("O613").Value2 = "Región del Maule"
$ws.Range("P613").Value2 = 389
$ws.Range("Q613").Value2 = 18
$ws.Range("R613").Value2 = "Hortaliza"
